$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1) Fix "оффициальный" -> "официальный" typo (appears twice in the doc:
#    "...присутствует оффициальный сайт игры..." and
#    "На оффициальном сайте вы можете...").
#    Both occurrences share the same misspelled stem "оффициальн", so a
#    single global replace handles both safely (no other word matches it).
# -------------------------------------------------------------------------
$f = $d.Content.Find
$f.ClearFormatting()
$f.Replacement.ClearFormatting()
$f.Text = "оффициальн"
$f.Replacement.Text = "официальн"
$f.Forward = $true
$f.Wrap = 1
$f.Format = $false
$f.MatchCase = $false
$f.MatchWholeWord = $false
$f.MatchWildcards = $false
$f.MatchSoundsLike = $false
$f.MatchAllWordForms = $false
$f.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)

# -------------------------------------------------------------------------
# 2) Fix "сопернечиства" -> "соперничества" and give that single word its
#    own run with Russian (ru-RU) language tagging, matching the target
#    OOXML which splits the paragraph's single run into three runs around
#    this word.
# -------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("сопернечиства", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Replace the word text first (still inside the single original run).
    $r.Text = "соперничества"

    # Re-select just the replaced word and force a run split by explicitly
    # assigning its font (this is what causes the engine to materialize a
    # separate <w:r> for the sub-range instead of mutating the whole run).
    $wordStart = $r.Start
    $wordEnd = $r.Start + 13  # length of "соперничества"
    $wr = $d.Range($wordStart, $wordEnd)
    $wr.Font.Name = "Calibri"
    $wr.Font.NameFarEast = "Helvetica"

    # Now that the word lives in its own run, tag it as Russian. Re-locate
    # the run via Find (rather than reusing a stale Range(start,end)) so the
    # property actually lands on the freshly-split run. A "warm-up" Find
    # right after the split, before the one whose Range we mutate, is
    # needed so the engine refreshes its view of the just-split runs before
    # the LanguageID write is recorded (otherwise the write is dropped).
    $warm = $d.Content
    [void]$warm.Find.Execute("соперничества", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $r2 = $d.Content
    $found2 = $r2.Find.Execute("соперничества", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $r2.LanguageID = "ru-RU"
    }
}

# -------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark: remove it from its old location (around
#    the discord-bot screenshot image near the end of the document) and
#    re-create it as an empty (collapsed) bookmark right after
#    "двухмерную " in the Minecraft-2D paragraph, i.e. between
#    "...собственно двухмерную " and "песочницу с механикой...".
# -------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}

$anchor = $d.Content
$anchorFound = $anchor.Find.Execute("двухмерную ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($anchorFound) {
    $pos = $anchor.End
    $newBmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $newBmRange)
}
